# Applies the "actualizado schema das classes e relatorio" edit:
#  - removes the two arrow-diagram labels that are no longer needed
#    ("setOnClickListener" in C9:D9, "MainActivity" in D11:E11)
#  - removes the now-unused box border around D11:E11
#  - removes the thick-bottom row formatting that Excel had reserved for
#    the border that used to sit under D11:E11
#  - deletes the vertical connector arrow that used to point at D9/D11
#  - leaves the cursor/selection on E3:F3 (where the author's session ended)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the two labels that were removed from the diagram.
$ws.Range("C9").Value = ""
$ws.Range("D11").Value = ""

# The merged D11:E11 cell no longer has a box border around it.
$ws.Range("D11:E11").Borders.LineStyle = -4142

# Row 10 no longer reserves extra height for a thick bottom border.
$ws.Rows.Item(10).AutoFit()

# The connector arrow that ran from D9 down to D11 is removed from the
# drawing along with the label it pointed to.
$ws.Shapes.Item("Straight Arrow Connector 4").Delete()

# Selection ends on E3:F3.
$ws.Range("E3:F3").Select()
